# Generate Report for Archive
# 1. Update the shared status string "Ready for handoff" -> "In Translation"
# 2. Shrink the (now narrower) status column widths accordingly

$wb = $excel.ActiveWorkbook

# --- 1. Replace the status text everywhere it appears (every sheet) ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Resize the status columns to fit the shorter text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E:F").ColumnWidth = 13.4101845877511

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C:C").ColumnWidth = 13.4101845877511

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C:C").ColumnWidth = 13.4101845877511
